$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2916266.02
$ws.Range("C7").Value = -34.36383648370334
$ws.Range("D7").Value = 2935
$ws.Range("E7").Value = 2935
$ws.Range("F7").Value = 993.6170425894378
$ws.Range("G7").Value = 5.912391963605113

$wb.Save()
